$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.676.21'
$ws.Range('E2').Value = '  -1.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.490.21'
$ws.Range('E3').Value = '  -1.13%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '534.37'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.24'
$ws.Range('E6').Value = '  -2.52%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.565'
$ws.Range('E8').Value = '  +0.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.511.12'
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.100'
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('E11').Value = '  -2.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.32'
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.345'
$ws.Range('E13').Value = '  -3.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.956.20'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.84'
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '58.669.65'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.506.16'
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.05'
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.26'
$ws.Range('E20').Value = '  -1.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '322.27'
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.88'
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.29'
$ws.Range('E24').Value = '  +3.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.419'
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.165'
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.53'
$ws.Range('E28').Value = '  -3.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.69'
$ws.Range('E29').Value = '  -3.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0763'
$ws.Range('E30').Value = '  -2.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.76'
$ws.Range('E31').Value = '  -1.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '166.93'
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('E33').Value = '  +3.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.46'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.37'
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.07'
$ws.Range('E37').Value = '  -4.78%  '
$ws.Range('E38').Value = '  -4.02%  '
$ws.Range('E39').Value = '  -0.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.810'
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.58'
$ws.Range('E41').Value = '  -2.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '284.41'
$ws.Range('E42').Value = '  +1.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.21'
$ws.Range('E43').Value = '  -0.94%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.995'
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '131.18'
$ws.Range('E45').Value = '  +6.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.602'
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.88'
$ws.Range('E47').Value = '  +0.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0922'
$ws.Range('E48').Value = '  -1.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0505'
$ws.Range('E49').Value = '  -1.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0219'
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.12'
$ws.Range('E51').Value = '  -4.12%  '
